$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new row (row 20) the same look as the rest of the table
# (rows 3-19 use a centered style) by copying the formatting from the
# row right above it.
$ws.Range("A19:C19").Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B20").Value = "逃离鸭科夫"
$ws.Range("C20").Value = 1276

# The "Date" column stores plain text like "2025/11/11", not a real date
# value. Assigning "2025/11/29" directly to .Value would let Excel's
# auto-detection turn it into a date serial number, which doesn't match
# the rest of the column. Instead, enter it as a string formula (so it is
# never parsed as a date) and then flatten the formula down to its static
# value in place.
$ws.Range("A20").Formula = "=""2025/11/29"""
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false
